$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: "Save" header, copy formatting from the neighboring header
# cell (G1) so it matches the other bold/bordered header cells, then set
# its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New column H data values
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
